$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5297.9
$ws.Range("I70").Value = 3424.75
$ws.Range("J70").Value = 6546.6665
$ws.Range("K70").Value = 10274.25
$ws.Range("L70").Value = 19639.9995
$ws.Range("M70").Value = -10004.25
$ws.Range("N70").Value = -20179.9995
$ws.Range("H73").Value = 5297.9
$ws.Range("I73").Value = 3424.75
$ws.Range("J73").Value = 6546.6665
$ws.Range("K73").Value = 10274.25
$ws.Range("L73").Value = 19639.9995
$ws.Range("M73").Value = -9338.25
$ws.Range("N73").Value = -21511.9995
$ws.Range("H74").Value = 7638.6665
$ws.Range("I74").Value = 6730.4614
$ws.Range("K74").Value = 6730.4614
$ws.Range("M74").Value = -5794.4614
$ws.Range("H77").Value = 7638.6665
$ws.Range("I77").Value = 6730.4614
$ws.Range("K77").Value = 33652.307
$ws.Range("M77").Value = -28972.307
$ws.Range("H100").Value = 19307.887
$ws.Range("I100").Value = 72755.86
$ws.Range("J100").Value = 9196.108
$ws.Range("K100").Value = 72755.86
$ws.Range("L100").Value = 9196.108
$ws.Range("M100").Value = -72214.86
$ws.Range("N100").Value = -10278.108
$ws.Range("H135").Value = 3852365.5
$ws.Range("I135").Value = 5556745.5
$ws.Range("J135").Value = 17510.625
$ws.Range("K135").Value = 50010709.5
$ws.Range("L135").Value = 157595.625
$ws.Range("M135").Value = -50008174.5
$ws.Range("N135").Value = -162665.625
$ws.Range("H137").Value = 7268.8335
$ws.Range("I137").Value = 4343.4375
$ws.Range("J137").Value = 9069.076999999999
$ws.Range("K137").Value = 13030.3125
$ws.Range("L137").Value = 27207.231
$ws.Range("M137").Value = -10480.3125
$ws.Range("N137").Value = -32307.231
$ws.Range("H138").Value = 4830.0938
$ws.Range("I138").Value = 4764.8335
$ws.Range("J138").Value = 4845.154
$ws.Range("K138").Value = 14294.5005
$ws.Range("L138").Value = 14535.462
$ws.Range("M138").Value = -9154.500499999998
$ws.Range("N138").Value = -24815.462

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2683440.5
$ws.Range("I32").Value = 2874392
$ws.Range("K32").Value = 2874392
$ws.Range("M32").Value = -2874105
$ws.Range("H74").Value = 205720.95
$ws.Range("I74").Value = 267477.25
$ws.Range("K74").Value = 267477.25
$ws.Range("M74").Value = -266603.25
$ws.Range("H77").Value = 205720.95
$ws.Range("I77").Value = 267477.25
$ws.Range("K77").Value = 1337386.25
$ws.Range("M77").Value = -1333018.25
$ws.Range("H97").Value = 9261009
$ws.Range("I97").Value = 12347679
$ws.Range("K97").Value = 12347679
$ws.Range("M97").Value = -12347183

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 50977.8
$ws.Range("J92").Value = 50977.8
$ws.Range("L92").Value = 50977.8
$ws.Range("N92").Value = -55969.8
$ws.Range("H94").Value = 1716.52
$ws.Range("I94").Value = 1431.4117
$ws.Range("J94").Value = 2322.375
$ws.Range("K94").Value = 1431.4117
$ws.Range("L94").Value = 2322.375
$ws.Range("M94").Value = -980.4117000000001
$ws.Range("N94").Value = -3224.375
$ws.Range("H99").Value = 1783.9286
$ws.Range("I99").Value = 1584.5555
$ws.Range("J99").Value = 2142.8
$ws.Range("K99").Value = 1584.5555
$ws.Range("L99").Value = 2142.8
$ws.Range("M99").Value = -86.55549999999994
$ws.Range("N99").Value = -5138.8
$ws.Range("H100").Value = 19749
$ws.Range("J100").Value = 19749
$ws.Range("L100").Value = 19749
$ws.Range("N100").Value = -21913

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3267.318
$ws.Range("I31").Value = 1727
$ws.Range("J31").Value = 4437.96
$ws.Range("K31").Value = 1727
$ws.Range("L31").Value = 4437.96
$ws.Range("M31").Value = -1432
$ws.Range("N31").Value = -5027.96
$ws.Range("H34").Value = 3267.318
$ws.Range("I34").Value = 1727
$ws.Range("J34").Value = 4437.96
$ws.Range("K34").Value = 1727
$ws.Range("L34").Value = 4437.96
$ws.Range("M34").Value = -1525
$ws.Range("N34").Value = -4841.96
$ws.Range("H62").Value = 26666.166
$ws.Range("I62").Value = 19998
$ws.Range("J62").Value = 27999.8
$ws.Range("K62").Value = 19998
$ws.Range("L62").Value = 27999.8
$ws.Range("M62").Value = -19374
$ws.Range("N62").Value = -29247.8
$ws.Range("H65").Value = 26666.166
$ws.Range("I65").Value = 19998
$ws.Range("J65").Value = 27999.8
$ws.Range("K65").Value = 99990
$ws.Range("L65").Value = 139999
$ws.Range("M65").Value = -96870
$ws.Range("N65").Value = -146239
$ws.Range("H92").Value = 35050.25
$ws.Range("J92").Value = 35050.25
$ws.Range("L92").Value = 35050.25
$ws.Range("N92").Value = -40042.25
$ws.Range("H134").Value = 6066.7407
$ws.Range("I134").Value = 5894.08
$ws.Range("J134").Value = 8225
$ws.Range("K134").Value = 17682.24
$ws.Range("L134").Value = 24675
$ws.Range("M134").Value = -15147.24
$ws.Range("N134").Value = -29745

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3562.2222
$ws.Range("I2").Value = 84.44444
$ws.Range("J2").Value = 5301.1113
$ws.Range("K2").Value = 506.66664
$ws.Range("L2").Value = 31806.6678
$ws.Range("M2").Value = -393.66664
$ws.Range("N2").Value = -32032.6678
$ws.Range("H6").Value = 1335.3334
$ws.Range("I6").Value = 1335.3334
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 4006.0002
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3893.0002
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 660.4286
$ws.Range("I9").Value = 660.4286
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1981.2858
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1757.2858
$ws.Range("N9").ClearContents()
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H22").Value = 563.3333
$ws.Range("I22").Value = 224.28572
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 672.85716
$ws.Range("L22").Value = 5250
$ws.Range("M22").Value = -503.85716
$ws.Range("N22").Value = -5588
$ws.Range("H27").Value = 563.3333
$ws.Range("I27").Value = 224.28572
$ws.Range("J27").Value = 1750
$ws.Range("K27").Value = 672.85716
$ws.Range("L27").Value = 5250
$ws.Range("M27").Value = -570.85716
$ws.Range("N27").Value = -5454
$ws.Range("H68").Value = 110367.69
$ws.Range("I68").Value = 252152.5
$ws.Range("J68").Value = 7251.4546
$ws.Range("K68").Value = 756457.5
$ws.Range("L68").Value = 21754.3638
$ws.Range("M68").Value = -755646.5
$ws.Range("N68").Value = -23376.3638
$ws.Range("H69").Value = 2824.75
$ws.Range("I69").Value = 900
$ws.Range("J69").Value = 4749.5
$ws.Range("K69").Value = 2700
$ws.Range("L69").Value = 14248.5
$ws.Range("M69").Value = -1889
$ws.Range("N69").Value = -15870.5
$ws.Range("H71").Value = 110367.69
$ws.Range("I71").Value = 252152.5
$ws.Range("J71").Value = 7251.4546
$ws.Range("K71").Value = 2269372.5
$ws.Range("L71").Value = 65263.0914
$ws.Range("M71").Value = -2265316.5
$ws.Range("N71").Value = -73375.0914
$ws.Range("H72").Value = 2824.75
$ws.Range("I72").Value = 900
$ws.Range("J72").Value = 4749.5
$ws.Range("K72").Value = 8100
$ws.Range("L72").Value = 42745.5
$ws.Range("M72").Value = -4044
$ws.Range("N72").Value = -50857.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2946.525
$ws.Range("I132").Value = 1859.1482
$ws.Range("J132").Value = 5204.923
$ws.Range("K132").Value = 5577.444600000001
$ws.Range("L132").Value = 15614.769
$ws.Range("M132").Value = -3047.444600000001
$ws.Range("N132").Value = -20674.769

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4667.4287
$ws.Range("I16").Value = 1334.4286
$ws.Range("J16").Value = 8000.4287
$ws.Range("K16").Value = 1334.4286
$ws.Range("L16").Value = 8000.4287
$ws.Range("M16").Value = -1164.4286
$ws.Range("N16").Value = -8340.4287
$ws.Range("H95").Value = 312000
$ws.Range("I95").Value = 312000
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 312000
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = -309254
$ws.Range("N95").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802
$ws.Range("H96").Value = 11279
$ws.Range("I96").Value = 2363.6667
$ws.Range("J96").Value = 21977.4
$ws.Range("K96").Value = 2363.6667
$ws.Range("L96").Value = 21977.4
$ws.Range("M96").Value = -990.6667000000002
$ws.Range("N96").Value = -24723.4
